$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.171470522880554
$ws.Range("B1").Value = 2.437279939651489
$ws.Range("D1").Value = 2.365739345550537
$ws.Range("E1").Value = 1.238474369049072
